$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 36.94436433333333
$ws.Range("H2").Value = 110.833093
$ws.Range("I2").Value = 0.8328964975864823
$ws.Range("J2").Value = 0.8328964975864824
$ws.Range("O2").Value = 0.1651823222987733
$ws.Range("P2").Value = 0.1651823222987733
$ws.Range("Q2").Value = 0.09684349370577777
$ws.Range("R2").Value = 0.871591443352
$ws.Range("S2").Value = 0.1375797777058498
$ws.Range("T2").Value = 0.1375797777058498

# Row 3
$ws.Range("G3").Value = 36.94436433333333
$ws.Range("H3").Value = 110.833093
$ws.Range("I3").Value = 0.8328964975864823
$ws.Range("J3").Value = 0.8328964975864824
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.013248
$ws.Range("N3").Value = 0.039744
$ws.Range("O3").Value = 0.8348176777012268
$ws.Range("P3").Value = 0.8348176777012267
$ws.Range("Q3").Value = 0.489438938688
$ws.Range("R3").Value = 4.404950448192
$ws.Range("S3").Value = 0.6953167198806327
$ws.Range("T3").Value = 0.6953167198806327

# Row 4
$ws.Range("I4").Value = 0.07608399754092349
$ws.Range("J4").Value = 0.07608399754092349
$ws.Range("O4").Value = 0.1651823222987733
$ws.Range("P4").Value = 0.1651823222987733
$ws.Range("S4").Value = 0.0125677314035839
$ws.Range("T4").Value = 0.0125677314035839

# Row 5
$ws.Range("I5").Value = 0.07608399754092349
$ws.Range("J5").Value = 0.07608399754092349
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.013248
$ws.Range("N5").Value = 0.039744
$ws.Range("O5").Value = 0.8348176777012268
$ws.Range("P5").Value = 0.8348176777012267
$ws.Range("Q5").Value = 0.044709602112
$ws.Range("R5").Value = 0.402386419008
$ws.Range("S5").Value = 0.0635162661373396
$ws.Range("T5").Value = 0.06351626613733959

# Row 6
$ws.Range("G6").Value = 4.037305666666668
$ws.Range("H6").Value = 12.111917
$ws.Range("I6").Value = 0.09101950487259411
$ws.Range("J6").Value = 0.09101950487259411
$ws.Range("O6").Value = 0.1651823222987733
$ws.Range("P6").Value = 0.1651823222987733
$ws.Range("Q6").Value = 0.01058312392088889
$ws.Range("R6").Value = 0.09524811528800001
$ws.Range("S6").Value = 0.01503481318933961
$ws.Range("T6").Value = 0.01503481318933961

# Row 7
$ws.Range("G7").Value = 4.037305666666668
$ws.Range("H7").Value = 12.111917
$ws.Range("I7").Value = 0.09101950487259411
$ws.Range("J7").Value = 0.09101950487259411
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.013248
$ws.Range("N7").Value = 0.039744
$ws.Range("O7").Value = 0.8348176777012268
$ws.Range("P7").Value = 0.8348176777012267
$ws.Range("Q7").Value = 0.05348622547200001
$ws.Range("R7").Value = 0.4813760292480001
$ws.Range("S7").Value = 0.07598469168325452
$ws.Range("T7").Value = 0.0759846916832545
